$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new literal text value, derived from the OOXML diff.
$changes = @(
    @{Cell='E2'; Value='2.57%'},
    @{Cell='G2'; Value='4'},
    @{Cell='D3'; Value='35.04'},
    @{Cell='E3'; Value='12.03%'},
    @{Cell='G3'; Value='4'},
    @{Cell='D4'; Value='5.142'},
    @{Cell='E4'; Value='4.73%'},
    @{Cell='G4'; Value='4'},
    @{Cell='E5'; Value='4.42%'},
    @{Cell='G5'; Value='4'},
    @{Cell='D6'; Value='2.337'},
    @{Cell='E6'; Value='4.17%'},
    @{Cell='G6'; Value='4'},
    @{Cell='D7'; Value='8.038'},
    @{Cell='E7'; Value='3.65%'},
    @{Cell='G7'; Value='4'},
    @{Cell='D8'; Value='3.949'},
    @{Cell='E8'; Value='5.45%'},
    @{Cell='G8'; Value='4'},
    @{Cell='D9'; Value='0.9295'},
    @{Cell='E9'; Value='1.56%'},
    @{Cell='G9'; Value='4'},
    @{Cell='D10'; Value='0.09981'},
    @{Cell='E10'; Value='11.45%'},
    @{Cell='G10'; Value='4'},
    @{Cell='D11'; Value='0.1784'},
    @{Cell='E11'; Value='3.92%'},
    @{Cell='G11'; Value='4'},
    @{Cell='D12'; Value='0.08625'},
    @{Cell='E12'; Value='3.69%'},
    @{Cell='G12'; Value='4'},
    @{Cell='D13'; Value='0.03326'},
    @{Cell='E13'; Value='6.90%'},
    @{Cell='G13'; Value='4'},
    @{Cell='D14'; Value='0.09918'},
    @{Cell='E14'; Value='-1.54%'},
    @{Cell='G14'; Value='4'},
    @{Cell='D15'; Value='0.001507'},
    @{Cell='E15'; Value='0.03%'},
    @{Cell='G15'; Value='4'},
    @{Cell='D16'; Value='0.005762'},
    @{Cell='E16'; Value='-0.40%'},
    @{Cell='G16'; Value='4'},
    @{Cell='D17'; Value='3.461'},
    @{Cell='E17'; Value='-1.30%'},
    @{Cell='G17'; Value='4'},
    @{Cell='D18'; Value='2.137'},
    @{Cell='E18'; Value='3.01%'},
    @{Cell='G18'; Value='4'},
    @{Cell='D19'; Value='0.3358'},
    @{Cell='E19'; Value='0.88%'},
    @{Cell='G19'; Value='4'},
    @{Cell='E20'; Value='2.74%'},
    @{Cell='G20'; Value='4'},
    @{Cell='D21'; Value='4.286'},
    @{Cell='E21'; Value='7.87%'},
    @{Cell='G21'; Value='4'},
    @{Cell='D22'; Value='0.2302'},
    @{Cell='E22'; Value='9.57%'},
    @{Cell='G22'; Value='4'},
    @{Cell='D23'; Value='0.04550'},
    @{Cell='E23'; Value='-0.15%'},
    @{Cell='G23'; Value='4'},
    @{Cell='D24'; Value='0.001215'},
    @{Cell='E24'; Value='0.11%'},
    @{Cell='G24'; Value='4'},
    @{Cell='D25'; Value='0.004377'},
    @{Cell='E25'; Value='-5.25%'},
    @{Cell='G25'; Value='4'},
    @{Cell='E26'; Value='0.01%'},
    @{Cell='G26'; Value='4'},
    @{Cell='E27'; Value='-0.04%'},
    @{Cell='G27'; Value='4'},
    @{Cell='G28'; Value='4'},
    @{Cell='G29'; Value='4'},
    @{Cell='G30'; Value='4'},
    @{Cell='G31'; Value='4'},
    @{Cell='G32'; Value='4'},
    @{Cell='G33'; Value='4'},
    @{Cell='G34'; Value='4'},
    @{Cell='G35'; Value='4'},
    @{Cell='G36'; Value='4'},
    @{Cell='G37'; Value='4'},
    @{Cell='G38'; Value='4'},
    @{Cell='D39'; Value='0.01791'},
    @{Cell='E39'; Value='10.90%'},
    @{Cell='G39'; Value='4'},
    @{Cell='D40'; Value='0.04795'},
    @{Cell='E40'; Value='6.92%'},
    @{Cell='G40'; Value='4'},
    @{Cell='D41'; Value='0.007802'},
    @{Cell='E41'; Value='6.59%'},
    @{Cell='G41'; Value='4'},
    @{Cell='D42'; Value='0.1410'},
    @{Cell='E42'; Value='6.05%'},
    @{Cell='G42'; Value='4'},
    @{Cell='D43'; Value='0.006913'},
    @{Cell='E43'; Value='-22.71%'},
    @{Cell='G43'; Value='4'},
    @{Cell='D44'; Value='0.002102'},
    @{Cell='E44'; Value='7.02%'},
    @{Cell='G44'; Value='4'},
    @{Cell='D45'; Value='0.009450'},
    @{Cell='E45'; Value='9.75%'},
    @{Cell='G45'; Value='4'},
    @{Cell='D46'; Value='0.00006117'},
    @{Cell='E46'; Value='1.28%'},
    @{Cell='G46'; Value='4'},
    @{Cell='E47'; Value='0.04%'},
    @{Cell='G47'; Value='4'},
    @{Cell='D48'; Value='3.037'},
    @{Cell='E48'; Value='31.78%'},
    @{Cell='G48'; Value='4'},
    @{Cell='D49'; Value='0.002002'},
    @{Cell='E49'; Value='0.05%'},
    @{Cell='G49'; Value='4'},
    @{Cell='E50'; Value='0.04%'},
    @{Cell='G50'; Value='4'},
    @{Cell='E51'; Value='0.04%'},
    @{Cell='G51'; Value='4'}
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    # Force the write to be stored as literal text (matches the source
    # inlineStr cells: e.g. "2.57%", "35.04", "0.001500", "4") instead of
    # letting Excel auto-coerce it into a Number/Percentage.
    $cell.NumberFormat = "@"
    $cell.Value = $change.Value
    # Drop the temporary Text number-format override again so the cell's
    # style stays the same as before the edit (no s="..." introduced).
    $cell.ClearFormats()
}
